# Edit script: clean up column headers, title-case Spanish connector words
# in state/municipality names, fix two floating point values, and remove
# the trailing metadata/footer rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename the header row (row 1) to the new machine-friendly names.
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value2 = "mx_state"
$ws.Cells.Item(1, 2).Value2 = "mx_municipality"
$ws.Cells.Item(1, 3).Value2 = "n_matriculas"
$ws.Cells.Item(1, 4).Value2 = "pct_matriculas"

# ---------------------------------------------------------------------
# 2) Title-case the Spanish connector words ("de", "del", "la", "los",
#    "las", "el", "y") inside the state (col A) and municipality (col B)
#    name cells, e.g. "Pabellón de Arteaga" -> "Pabellón De Arteaga".
# ---------------------------------------------------------------------
$lastDataRow = 939

for ($r = 2; $r -le $lastDataRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $aVal = $aCell.Value2
    if ($aVal -ne $null) {
        $s = $aVal
        $s = $s -creplace '\bde\b', 'De'
        $s = $s -creplace '\bdel\b', 'Del'
        $s = $s -creplace '\bla\b', 'La'
        $s = $s -creplace '\blos\b', 'Los'
        $s = $s -creplace '\blas\b', 'Las'
        $s = $s -creplace '\bel\b', 'El'
        $s = $s -creplace '\by\b', 'Y'
        $aCell.Value2 = $s
    }

    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value2
    if ($bVal -ne $null) {
        $s = $bVal
        $s = $s -creplace '\bde\b', 'De'
        $s = $s -creplace '\bdel\b', 'Del'
        $s = $s -creplace '\bla\b', 'La'
        $s = $s -creplace '\blos\b', 'Los'
        $s = $s -creplace '\blas\b', 'Las'
        $s = $s -creplace '\bel\b', 'El'
        $s = $s -creplace '\by\b', 'Y'
        $bCell.Value2 = $s
    }
}

# ---------------------------------------------------------------------
# 3) Tiny floating point re-expressions for two cells (last-bit rounding
#    differences produced by the original recalculation).
# ---------------------------------------------------------------------
$ws.Cells.Item(98, 4).Value2 = 0.009978617248752672
$ws.Cells.Item(99, 4).Value2 = 0.009087669280114039

# ---------------------------------------------------------------------
# 4) Remove the trailing metadata/footer rows (941-945), shrinking the
#    sheet's dimension from A1:D945 down to A1:D939.
# ---------------------------------------------------------------------
$ws.Range("A941:A945").EntireRow.Delete()
